$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting rows 59:137 down to 60:138.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record's data.
$ws.Range("A59").Value = 8
$ws.Range("B59").Value = "Terminal La Palmera de La Serena"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44665
$ws.Range("D59").NumberFormat = $ws.Range("D58").NumberFormat
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 100112044
$ws.Range("G59").Value = "Perejil"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 3000
$ws.Range("K59").Value = 2000
$ws.Range("L59").Value = 2500
$ws.Range("M59").Value = 2250
$ws.Range("N59").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O59").Value = "Provincia del Elquí"
$ws.Range("P59").Value = 1500
$ws.Range("Q59").Value = 1.5
$ws.Range("R59").Value = "Hortaliza"
